# ReportingOrganisationGroup.xlsx edit:
# The codeforiati:group-name / codeforiati:group-code columns (D and E) were
# swapped for every data row (including the header row), so that the
# "group-code" column now precedes the "group-name" column's values in the
# underlying shared-string ordering. In terms of the rendered worksheet this
# is simply: for every used row, swap the contents of column D and column E.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)

    $dVal = $dCell.Value2
    $eVal = $eCell.Value2

    $dCell.Value2 = $eVal
    $eCell.Value2 = $dVal
}
